$d = $word.ActiveDocument
$sec = $d.Sections.First

# The document has a "different first page" header/footer setup, so each
# section exposes two in-use HeaderFooter stories:
#   Headers.Item(1) / Footers.Item(1) -> the "default" (non-first-page) story
#   Headers.Item(2) / Footers.Item(2) -> the "first page" story
# Each of the four stories holds exactly one inline picture, and every one
# of those pictures needs its display name updated:
#   - the BTec logo pictures (in the headers) go from "image1.jpg" to "image2.jpg"
#   - the Pearson logo pictures (in the footers) go from "image2.png" to "image1.png"

$header1 = $sec.Headers.Item(1)
if ($header1.Exists) {
    $header1.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

$header2 = $sec.Headers.Item(2)
if ($header2.Exists) {
    $header2.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

$footer1 = $sec.Footers.Item(1)
if ($footer1.Exists) {
    $footer1.Range.InlineShapes.Item(1).Name = "image1.png"
}

$footer2 = $sec.Footers.Item(2)
if ($footer2.Exists) {
    $footer2.Range.InlineShapes.Item(1).Name = "image1.png"
}
